$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-17 07:48:41'
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '45%'
$ws.Range("M2").Value = '2.0 °C 7:24 TU'
$ws.Range("O2").Value = '0.7 °C'
$ws.Range("E3").Value = '2026-02-17 07:48:44'
$ws.Range("E4").Value = '2026-02-17 07:48:46'
$ws.Range("J4").Value = '1016.2 hPa'
$ws.Range("K4").Value = '0.0 MJ/m2'
$ws.Range("O4").Value = '6.0 °C'
$ws.Range("E5").Value = '2026-02-17 07:48:49'
$ws.Range("M5").Value = '-5.0 °C 7:11 TU'
$ws.Range("E6").Value = '2026-02-17 07:48:51'
$ws.Range("J6").Value = '1015.9 hPa'
$ws.Range("K6").Value = '0.1 MJ/m2'
$ws.Range("E7").Value = '2026-02-17 07:48:54'
$ws.Range("J7").Value = '1015.5 hPa'
$ws.Range("K7").Value = '0.0 MJ/m2'
$ws.Range("N7").Value = '12.1 °C 7:27 TU'
$ws.Range("O7").Value = '13.7 °C'
$ws.Range("E8").Value = '2026-02-17 07:48:56'
$ws.Range("J8").Value = '1015.7 hPa'
$ws.Range("K8").Value = '0.1 MJ/m2'
$ws.Range("N8").Value = '8.2 °C 7:23 TU'
$ws.Range("O8").Value = '9.6 °C'
$ws.Range("E9").Value = '2026-02-17 07:48:58'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '49%'
$ws.Range("K9").Value = '0.1 MJ/m2'
$ws.Range("O9").Value = '11.9 °C'
$ws.Range("E10").Value = '2026-02-17 07:49:01'
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = '82%'
$ws.Range("K10").Value = '0.1 MJ/m2'
$ws.Range("L10").Value = '19.8 km/h - 45º 7:11 TU'
$ws.Range("O10").Value = '8.3 °C'
$ws.Range("E11").Value = '2026-02-17 07:49:04'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '43%'
$ws.Range("O11").Value = '5.2 °C'
$ws.Range("E12").Value = '2026-02-17 07:49:06'
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '51%'
$ws.Range("O12").Value = '12.2 °C'
$ws.Range("E13").Value = '2026-02-17 07:49:09'
$ws.Range("J13").Value = '1017.4 hPa'
$ws.Range("K13").Value = '0.0 MJ/m2'
$ws.Range("O13").Value = '4.1 °C'
$ws.Range("E14").Value = '2026-02-17 07:49:11'
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = '63%'
$ws.Range("N14").Value = '7.8 °C 7:19 TU'
$ws.Range("O14").Value = '12.2 °C'
$ws.Range("E15").Value = '2026-02-17 07:49:14'
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '51%'
$ws.Range("N15").Value = '9.8 °C 7:26 TU'
$ws.Range("O15").Value = '11.7 °C'
$ws.Range("E16").Value = '2026-02-17 07:49:16'
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '52%'
$ws.Range("M16").Value = '-4.2 °C 7:24 TU'
$ws.Range("E17").Value = '2026-02-17 07:49:18'
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = '53%'
$ws.Range("N17").Value = '-0.3 °C 7:27 TU'
$ws.Range("O17").Value = '2.3 °C'
$ws.Range("E18").Value = '2026-02-17 07:49:21'
$ws.Range("J18").Value = '1016.2 hPa'
$ws.Range("K18").Value = '0.1 MJ/m2'
$ws.Range("E19").Value = '2026-02-17 07:49:23'
$ws.Range("K19").Value = '0.0 MJ/m2'
$ws.Range("O19").Value = '5.7 °C'
$ws.Range("E20").Value = '2026-02-17 07:49:25'
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '40%'
$ws.Range("K20").Value = '0.1 MJ/m2'
$ws.Range("E21").Value = '2026-02-17 07:49:27'
$ws.Range("J21").Value = '1016.3 hPa'
$ws.Range("K21").Value = '0.0 MJ/m2'
$ws.Range("E22").Value = '2026-02-17 07:49:30'
$ws.Range("E23").Value = '2026-02-17 07:49:33'
$ws.Range("M23").Value = '-5.4 °C 7:29 TU'
$ws.Range("E24").Value = '2026-02-17 07:49:35'
$ws.Range("J24").Value = '1017.8 hPa'
$ws.Range("N24").Value = '9.1 °C 7:12 TU'
$ws.Range("O24").Value = '9.7 °C'
$ws.Range("E25").Value = '2026-02-17 07:49:38'
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = '43%'
$ws.Range("E26").Value = '2026-02-17 07:49:41'
$ws.Range("E27").Value = '2026-02-17 07:49:43'
$ws.Range("E28").Value = '2026-02-17 07:49:46'
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '91%'
$ws.Range("J28").Value = '1016.4 hPa'
$ws.Range("O28").Value = '5.2 °C'
$ws.Range("E29").Value = '2026-02-17 07:49:48'
$ws.Range("O29").Value = '11.4 °C'
$ws.Range("E30").Value = '2026-02-17 07:49:51'
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '51%'
$ws.Range("J30").Value = '1015.5 hPa'
$ws.Range("K30").Value = '0.1 MJ/m2'
$ws.Range("O30").Value = '11.3 °C'
$ws.Range("E31").Value = '2026-02-17 07:49:53'
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '67%'
$ws.Range("J31").Value = '1016.0 hPa'
$ws.Range("O31").Value = '9.4 °C'
$ws.Range("E32").Value = '2026-02-17 07:49:56'
$ws.Range("O32").Value = '6.2 °C'
$ws.Range("E33").Value = '2026-02-17 07:49:58'
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = '38%'
$ws.Range("J33").Value = '1016.4 hPa'
$ws.Range("K33").Value = '0.0 MJ/m2'
$ws.Range("E34").Value = '2026-02-17 07:50:01'
$ws.Range("M34").Value = '0.8 °C 7:24 TU'
$ws.Range("E35").Value = '2026-02-17 07:50:04'
$ws.Range("J35").Value = '1018.7 hPa'
$ws.Range("N35").Value = '4.8 °C 7:03 TU'
$ws.Range("E36").Value = '2026-02-17 07:50:06'
$ws.Range("J36").Value = '1015.9 hPa'
$ws.Range("K36").Value = '0.1 MJ/m2'
$ws.Range("O36").Value = '12.3 °C'
$ws.Range("E37").Value = '2026-02-17 07:50:09'
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = '58%'
$ws.Range("J37").Value = '1016.6 hPa'
$ws.Range("O37").Value = '7.1 °C'
$ws.Range("E38").Value = '2026-02-17 07:50:12'
$ws.Range("E39").Value = '2026-02-17 07:50:14'
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = '50%'
$ws.Range("K39").Value = '0.1 MJ/m2'
$ws.Range("M39").Value = '-2.5 °C 7:19 TU'
$ws.Range("O39").Value = '-4.4 °C'
$ws.Range("E40").Value = '2026-02-17 07:50:17'
$ws.Range("J40").Value = '1017.9 hPa'
$ws.Range("E41").Value = '2026-02-17 07:50:19'
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = '52%'
$ws.Range("J41").Value = '1016.1 hPa'
$ws.Range("K41").Value = '0.0 MJ/m2'
$ws.Range("O41").Value = '14.3 °C'
$ws.Range("E42").Value = '2026-02-17 07:50:22'
$ws.Range("O42").Value = '12.3 °C'
$ws.Range("E43").Value = '2026-02-17 07:50:24'
$ws.Range("E44").Value = '2026-02-17 07:50:27'
$ws.Range("E45").Value = '2026-02-17 07:50:30'
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = '54%'
$ws.Range("J45").Value = '1021.2 hPa'
$ws.Range("O45").Value = '4.4 °C'
$ws.Range("E46").Value = '2026-02-17 07:50:33'
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = '60%'
$ws.Range("J46").Value = '1018.1 hPa'
$ws.Range("K46").Value = '0.0 MJ/m2'
$ws.Range("N46").Value = '12.5 °C 7:00 TU'
